$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2359249329758713
$ws.Range("C2").Value = 0.5013404825737265
$ws.Range("J2").Value = 0.02144772117962467
$ws.Range("P2").Value = 0.1554959785522788
$ws.Range("S2").Value = 0.08579088471849866
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.005263157894736842
$ws.Range("J3").Value = 0.01578947368421053
$ws.Range("P3").Value = 0.7157894736842105
$ws.Range("S3").Value = 0.2526315789473684
$ws.Range("J4").Value = 0.07017543859649122
$ws.Range("O4").Value = 0.01754385964912281
$ws.Range("S4").Value = 0.2456140350877193
$ws.Range("B6").Value = 0.06060606060606061
$ws.Range("D6").Value = 0.0202020202020202
$ws.Range("E6").Value = 0.005050505050505051
$ws.Range("F6").Value = 0.08080808080808081
$ws.Range("J6").Value = 0.2828282828282828
$ws.Range("O6").Value = 0.0101010101010101
$ws.Range("Q6").Value = 0.1868686868686869
$ws.Range("R6").Value = 0.0505050505050505
$ws.Range("S6").Value = 0.303030303030303
$ws.Range("B7").Value = 0.1
$ws.Range("F7").Value = 0.07333333333333333
$ws.Range("J7").Value = 0.1733333333333333
$ws.Range("O7").Value = 0.02
$ws.Range("Q7").Value = 0.1333333333333333
$ws.Range("R7").Value = 0.07333333333333333
$ws.Range("S7").Value = 0.4266666666666667
$ws.Range("B8").Value = 0.1312910284463895
$ws.Range("D8").Value = 0.03063457330415755
$ws.Range("E8").Value = 0.00437636761487965
$ws.Range("F8").Value = 0.04595185995623632
$ws.Range("J8").Value = 0.1422319474835886
$ws.Range("O8").Value = 0.02844638949671772
$ws.Range("Q8").Value = 0.1553610503282276
$ws.Range("R8").Value = 0.07439824945295405
$ws.Range("S8").Value = 0.387308533916849
$ws.Range("B9").Value = 0.0738255033557047
$ws.Range("D9").Value = 0.02684563758389262
$ws.Range("F9").Value = 0.02013422818791946
$ws.Range("J9").Value = 0.1677852348993289
$ws.Range("Q9").Value = 0.2416107382550336
$ws.Range("R9").Value = 0.1006711409395973
$ws.Range("S9").Value = 0.3691275167785235
$ws.Range("B10").Value = 0.1330882352941176
$ws.Range("D10").Value = 0.025
$ws.Range("E10").Value = 0.002205882352941176
$ws.Range("F10").Value = 0.05882352941176471
$ws.Range("J10").Value = 0.1154411764705882
$ws.Range("O10").Value = 0.02352941176470588
$ws.Range("Q10").Value = 0.2161764705882353
$ws.Range("R10").Value = 0.09117647058823529
$ws.Range("S10").Value = 0.3345588235294117
$ws.Range("G11").Value = 0.1526104417670683
$ws.Range("J11").Value = 0.08835341365461848
$ws.Range("K11").Value = 0.2088353413654618
$ws.Range("L11").Value = 0.5341365461847389
$ws.Range("S11").Value = 0.01606425702811245
$ws.Range("G12").Value = 0.7185185185185186
$ws.Range("J12").Value = 0.2074074074074074
$ws.Range("L12").Value = 0.01481481481481482
$ws.Range("S12").Value = 0.05925925925925926
$ws.Range("G13").Value = 0.6285714285714286
$ws.Range("J13").Value = 0.3142857142857143
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.01851851851851852
$ws.Range("H15").Value = 0.1342592592592593
$ws.Range("I15").Value = 0.03703703703703703
$ws.Range("J15").Value = 0.4398148148148148
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.01851851851851852
$ws.Range("O15").Value = 0.01851851851851852
$ws.Range("S15").Value = 0.2777777777777778
$ws.Range("F16").Value = 0.01769911504424779
$ws.Range("H16").Value = 0.168141592920354
$ws.Range("I16").Value = 0.06637168141592921
$ws.Range("J16").Value = 0.4424778761061947
$ws.Range("K16").Value = 0.05752212389380531
$ws.Range("M16").Value = 0.02212389380530973
$ws.Range("N16").Value = 0.004424778761061947
$ws.Range("O16").Value = 0.07964601769911504
$ws.Range("S16").Value = 0.1415929203539823
$ws.Range("F17").Value = 0.01098901098901099
$ws.Range("H17").Value = 0.2065934065934066
$ws.Range("I17").Value = 0.08131868131868132
$ws.Range("J17").Value = 0.4527472527472527
$ws.Range("K17").Value = 0.06813186813186813
$ws.Range("M17").Value = 0.01978021978021978
$ws.Range("N17").Value = 0.002197802197802198
$ws.Range("O17").Value = 0.06373626373626373
$ws.Range("S17").Value = 0.0945054945054945
$ws.Range("F18").Value = 0.02590673575129534
$ws.Range("H18").Value = 0.1450777202072539
$ws.Range("I18").Value = 0.07253886010362694
$ws.Range("J18").Value = 0.5077720207253886
$ws.Range("K18").Value = 0.07253886010362694
$ws.Range("M18").Value = 0.0155440414507772
$ws.Range("O18").Value = 0.05699481865284974
$ws.Range("S18").Value = 0.1036269430051813
$ws.Range("F19").Value = 0.02040816326530612
$ws.Range("H19").Value = 0.2346938775510204
$ws.Range("I19").Value = 0.06462585034013606
$ws.Range("J19").Value = 0.3937074829931973
$ws.Range("K19").Value = 0.09948979591836735
$ws.Range("M19").Value = 0.01530612244897959
$ws.Range("N19").Value = 0.001700680272108843
$ws.Range("O19").Value = 0.06802721088435375
$ws.Range("S19").Value = 0.1020408163265306
